$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 93, shifting existing rows 93:105 down to 94:106.
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with the new data record.
$ws.Cells.Item(93, 1).Value = 9
$ws.Cells.Item(93, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(93, 3).Value = "Metropolitana"
$ws.Cells.Item(93, 4).Value = 44476
$ws.Cells.Item(93, 5).Value = 13
$ws.Cells.Item(93, 6).Value = "Fruta"
$ws.Cells.Item(93, 7).Value = 100101
$ws.Cells.Item(93, 8).Value = "Berries"
$ws.Cells.Item(93, 9).Value = 100101001
$ws.Cells.Item(93, 10).Value = "Arándano (blue)"
$ws.Cells.Item(93, 11).Value = "Sin especificar"
$ws.Cells.Item(93, 12).Value = "Primera"
$ws.Cells.Item(93, 13).Value = 300
$ws.Cells.Item(93, 14).Value = 14000
$ws.Cells.Item(93, 15).Value = 14000
$ws.Cells.Item(93, 16).Value = 14000
$ws.Cells.Item(93, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(93, 18).Value = "Provincia de Linares"
$ws.Cells.Item(93, 19).Value = 7000
$ws.Cells.Item(93, 20).Value = 2
